$wb = $excel.ActiveWorkbook

# --- Overview sheet: bump "Latest HO Xliff Generate Date" for the
#     "Ready for handoff" rows (4-7) from 16:47:18 -> 16:47:35 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4:G7").Value = "2016-09-07 16:47:35"

# --- zh-cn sheet: Priority low -> ht, Latest Handoff Datetime bump ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4:E7").Value = "ht"
$wsZhCn.Range("H4:H7").Value = "2016-09-07 16:47:29"

# --- de-de sheet: Priority low -> ht, Latest Handoff Datetime bump ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4:E7").Value = "ht"
$wsDeDe.Range("H4:H7").Value = "2016-09-07 16:47:35"
